$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) from 45184 -> 45186 for all data rows (2..483)
for ($r = 2; $r -le 483; $r++) {
    $ws.Range("C$r").Value = 45186
}

# 2. Add the display-text second argument to every HYPERLINK formula in rows 2..13
#    (columns S, T, U, V, W, X, Y) using the row's Beteckning (column A) value.
$hyperlinkCols = @("S","T","U","V","W","X","Y")
for ($r = 2; $r -le 13; $r++) {
    $idVal = $ws.Range("A$r").Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$r")
        $oldFormula = $cell.Formula
        if ($oldFormula -ne $null -and $oldFormula -ne "") {
            $newFormula = $oldFormula.Substring(0, $oldFormula.Length - 1) + ', "' + $idVal + '")'
            $cell.Formula = $newFormula
        }
    }
}

# 3. Row 483 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(483).RowHeight = 15

# 4. Append two new data rows (484 and 485)

# Row 484
$ws.Range("A484").Value = "A 43489-2023"
$ws.Range("B484").Value = 45184
$ws.Range("B484").NumberFormat = "YYYY-MM-DD"
$ws.Range("C484").Value = 45186
$ws.Range("C484").NumberFormat = "YYYY-MM-DD"
$ws.Range("D484").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E484").Value = "EKSJÖ"
$ws.Range("G484").Value = 1.4
$ws.Range("H484").Value = 0
$ws.Range("I484").Value = 0
$ws.Range("J484").Value = 0
$ws.Range("K484").Value = 0
$ws.Range("L484").Value = 0
$ws.Range("M484").Value = 0
$ws.Range("N484").Value = 0
$ws.Range("O484").Value = 0
$ws.Range("P484").Value = 0
$ws.Range("Q484").Value = 0
$ws.Range("R484").Value = ""
$ws.Range("R484").WrapText = $true
$ws.Rows.Item(484).RowHeight = 15

# Row 485
$ws.Range("A485").Value = "A 43483-2023"
$ws.Range("B485").Value = 45184
$ws.Range("B485").NumberFormat = "YYYY-MM-DD"
$ws.Range("C485").Value = 45186
$ws.Range("C485").NumberFormat = "YYYY-MM-DD"
$ws.Range("D485").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E485").Value = "EKSJÖ"
$ws.Range("G485").Value = 1.2
$ws.Range("H485").Value = 0
$ws.Range("I485").Value = 0
$ws.Range("J485").Value = 0
$ws.Range("K485").Value = 0
$ws.Range("L485").Value = 0
$ws.Range("M485").Value = 0
$ws.Range("N485").Value = 0
$ws.Range("O485").Value = 0
$ws.Range("P485").Value = 0
$ws.Range("Q485").Value = 0
$ws.Range("R485").Value = ""
$ws.Range("R485").WrapText = $true

Write-Host "Edit complete"
